# #fix: exclude troca pe cv preventivo
# Adjust quantity (F) and value (H) for rows where "troca pe cv preventivo"
# movements must be excluded from the stock totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 46;  F = 207.000;    H = 115666.50 },
    @{ Row = 56;  F = 1068.000;   H = 9731.75 },
    @{ Row = 58;  F = 2936.000;   H = 30249.47 },
    @{ Row = 70;  F = 821.000;    H = 4267.74 },
    @{ Row = 79;  F = 91.000;     H = 1690.75 },
    @{ Row = 82;  F = 1482.000;   H = 9193.41 },
    @{ Row = 83;  F = 2716.000;   H = 70744.33 },
    @{ Row = 97;  F = 2366.000;   H = 325904.67 },
    @{ Row = 113; F = 532.000;    H = 19742.96 },
    @{ Row = 115; F = 952.000;    H = 31365.50 },
    @{ Row = 120; F = 6599.000;   H = 25511.61 },
    @{ Row = 132; F = 251.000;    H = 1975.18 },
    @{ Row = 155; F = 497.000;    H = 2199.90 },
    @{ Row = 161; F = 54.000;     H = 7209.48 },
    @{ Row = 175; F = 2444.000;   H = 45763.90 },
    @{ Row = 178; F = 30.000;     H = 5740.16 },
    @{ Row = 179; F = 668.000;    H = 116644.22 },
    @{ Row = 198; F = 257.000;    H = 17795.97 },
    @{ Row = 233; F = 21615.800;  H = 49061.00 },
    @{ Row = 237; F = 564.600;    H = 5760.98 },
    @{ Row = 249; F = 4538.000;   H = 11391.96 }
)

foreach ($change in $changes) {
    $r = $change.Row
    $ws.Cells.Item($r, 6).Value = $change.F
    $ws.Cells.Item($r, 8).Value = $change.H
}
